$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1140.125
$ws.Range("I41").Value = 1925
$ws.Range("J41").Value = 878.5
$ws.Range("K41").Value = 1925
$ws.Range("L41").Value = 878.5
$ws.Range("M41").Value = -1485
$ws.Range("N41").Value = -1758.5

$ws.Range("H86").Value = 2794.9092
$ws.Range("I86").Value = 2078
$ws.Range("J86").Value = 3392.3333
$ws.Range("K86").Value = 2078
$ws.Range("L86").Value = 3392.3333
$ws.Range("M86").Value = -955
$ws.Range("N86").Value = -5638.3333

$ws.Range("H89").Value = 2794.9092
$ws.Range("I89").Value = 2078
$ws.Range("J89").Value = 3392.3333
$ws.Range("K89").Value = 10390
$ws.Range("L89").Value = 16961.6665
$ws.Range("M89").Value = -4774
$ws.Range("N89").Value = -28193.6665

$ws.Range("H92").Value = 934.2778
$ws.Range("J92").Value = 841.5833
$ws.Range("L92").Value = 841.5833
$ws.Range("N92").Value = -3337.5833

$ws.Range("H96").Value = 577.1111
$ws.Range("I96").Value = 705.2857
$ws.Range("K96").Value = 2115.8571
$ws.Range("M96").Value = -742.8571000000002

$ws.Range("H132").Value = 18875.117
$ws.Range("I132").Value = 3163.08
$ws.Range("J132").Value = 62519.668
$ws.Range("K132").Value = 9489.24
$ws.Range("L132").Value = 187559.004
$ws.Range("M132").Value = -6959.24
$ws.Range("N132").Value = -192619.004

$ws.Range("H135").Value = 29413928
$ws.Range("I135").Value = 33335252
$ws.Range("J135").Value = 3999.5
$ws.Range("K135").Value = 300017268
$ws.Range("L135").Value = 35995.5
$ws.Range("M135").Value = -300014733
$ws.Range("N135").Value = -41065.5

$ws.Range("H138").Value = 3839.549
$ws.Range("J138").Value = 5645.8887
$ws.Range("L138").Value = 16937.6661
$ws.Range("N138").Value = -27217.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16397810
$ws.Range("I32").Value = 16953604
$ws.Range("K32").Value = 16953604
$ws.Range("M32").Value = -16953317

$ws.Range("H45").Value = 2676.5715
$ws.Range("I45").Value = 2147.4
$ws.Range("K45").Value = 2147.4
$ws.Range("M45").Value = -1770.4

$ws.Range("H74").Value = 2778.077
$ws.Range("I74").Value = 2778.077
$ws.Range("K74").Value = 2778.077
$ws.Range("M74").Value = -1904.077

$ws.Range("H77").Value = 2778.077
$ws.Range("I77").Value = 2778.077
$ws.Range("K77").Value = 13890.385
$ws.Range("M77").Value = -9522.385000000002

$ws.Range("H102").Value = 10492.111
$ws.Range("I102").Value = 10427.25
$ws.Range("K102").Value = 10427.25
$ws.Range("M102").Value = -8805.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3000.25
$ws.Range("I20").Value = 917
$ws.Range("J20").Value = 5916.8
$ws.Range("K20").Value = 917
$ws.Range("L20").Value = 5916.8
$ws.Range("M20").Value = -670
$ws.Range("N20").Value = -6410.8

$ws.Range("H99").Value = 39591
$ws.Range("I99").Value = 47658.11
$ws.Range("J99").Value = 3289
$ws.Range("K99").Value = 47658.11
$ws.Range("L99").Value = 3289
$ws.Range("M99").Value = -46160.11
$ws.Range("N99").Value = -6285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1726.2903
$ws.Range("I31").Value = 1717.1666
$ws.Range("K31").Value = 1717.1666
$ws.Range("M31").Value = -1422.1666

$ws.Range("H34").Value = 1726.2903
$ws.Range("I34").Value = 1717.1666
$ws.Range("K34").Value = 1717.1666
$ws.Range("M34").Value = -1515.1666

$ws.Range("H94").Value = 12649
$ws.Range("I94").Value = 33411.668
$ws.Range("J94").Value = 2267.6667
$ws.Range("K94").Value = 33411.668
$ws.Range("L94").Value = 2267.6667
$ws.Range("M94").Value = -32960.668
$ws.Range("N94").Value = -3169.6667

$ws.Range("H105").Value = 2280.4
$ws.Range("I105").Value = 2367.4167
$ws.Range("J105").Value = 1932.3334
$ws.Range("K105").Value = 2367.4167
$ws.Range("L105").Value = 1932.3334
$ws.Range("M105").Value = -620.4167000000002
$ws.Range("N105").Value = -5426.3334

$ws.Range("H132").Value = 2948
$ws.Range("I132").Value = 2790.1538
$ws.Range("K132").Value = 8370.4614
$ws.Range("M132").Value = -5840.4614

$ws.Range("H134").Value = 2465.0908
$ws.Range("I134").Value = 2041.6875
$ws.Range("J134").Value = 16014
$ws.Range("K134").Value = 6125.0625
$ws.Range("L134").Value = 48042
$ws.Range("M134").Value = -3590.0625
$ws.Range("N134").Value = -53112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3700.5
$ws.Range("J75").Value = 5000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16996

$ws.Range("H78").Value = 3700.5
$ws.Range("J78").Value = 5000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -54984

$ws.Range("H95").Value = 7500
$ws.Range("J95").Value = 7500
$ws.Range("L95").Value = 22500
$ws.Range("N95").Value = -26618

$ws.Range("H137").Value = 2612.5
$ws.Range("I137").Value = 2283.9285
$ws.Range("K137").Value = 6851.7855
$ws.Range("M137").Value = -1751.7855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 957.375
$ws.Range("I16").Value = 714.5454999999999
$ws.Range("K16").Value = 714.5454999999999
$ws.Range("M16").Value = -544.5454999999999

$ws.Range("H93").Value = 7064.2856
$ws.Range("I93").Value = 6756.2144
$ws.Range("J93").Value = 7680.4287
$ws.Range("K93").Value = 6756.2144
$ws.Range("L93").Value = 7680.4287
$ws.Range("M93").Value = -5508.2144
$ws.Range("N93").Value = -10176.4287

$ws.Range("H132").Value = 5986.857
$ws.Range("I132").Value = 2818
$ws.Range("K132").Value = 8454
$ws.Range("M132").Value = -5924

$ws.Range("H136").Value = 4276.385
$ws.Range("I136").Value = 4276.385
$ws.Range("K136").Value = 12829.155
$ws.Range("M136").Value = -10279.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10172.777
$ws.Range("J41").Value = 10444.375
$ws.Range("L41").Value = 10444.375
$ws.Range("N41").Value = -11224.375

$ws.Range("H132").Value = 2899.2
$ws.Range("J132").Value = 2498.5
$ws.Range("L132").Value = 7495.5
$ws.Range("N132").Value = -12555.5

$ws.Range("H135").Value = 71042.625
$ws.Range("J135").Value = 71042.625
$ws.Range("L135").Value = 71042.625
$ws.Range("N135").Value = -81182.625

$ws.Range("H136").Value = 1777.7858
$ws.Range("I136").Value = 1478.3334
$ws.Range("K136").Value = 4435.0002
$ws.Range("M136").Value = -1885.0002
